# Update the "想去人数" (want-to-go count) figures on the "展览" and
# "全部类型" sheets to reflect newly generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1592
    $ws.Range("F3").Value = 130
    $ws.Range("F4").Value = 76
}
